$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44999
$ws.Range("J2").Value = 1100
$ws.Range("D3").Value = 45070
$ws.Range("J3").Value = 800
$ws.Range("D4").Value = 44881
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 1900
$ws.Range("M4").Value = 1950
$ws.Range("P4").Value = 650
$ws.Range("D5").Value = 44965
$ws.Range("J5").Value = 1120
$ws.Range("D6").Value = 45077
$ws.Range("J6").Value = 760
$ws.Range("D7").Value = 45091
$ws.Range("J7").Value = 800
$ws.Range("D8").Value = 45084
$ws.Range("J8").Value = 900
$ws.Range("D9").Value = 44985
$ws.Range("D10").Value = 45041
$ws.Range("J10").Value = 1160
$ws.Range("K10").Value = 2000
$ws.Range("L10").Value = 2500
$ws.Range("M10").Value = 2250
$ws.Range("P10").Value = 750
$ws.Range("D11").Value = 44953
$ws.Range("J11").Value = 1000
$ws.Range("D12").Value = 44964
$ws.Range("J12").Value = 1000
$ws.Range("D13").Value = 44883
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 1800
$ws.Range("M13").Value = 1900
$ws.Range("P13").Value = 633
$ws.Range("D14").Value = 45062
$ws.Range("J14").Value = 1100
$ws.Range("K14").Value = 2000
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = 2250
$ws.Range("P14").Value = 750
$ws.Range("D15").Value = 44971
$ws.Range("J15").Value = 1000
$ws.Range("K15").Value = 2000
$ws.Range("L15").Value = 2500
$ws.Range("M15").Value = 2250
$ws.Range("P15").Value = 750
$ws.Range("D16").Value = 44951
$ws.Range("J16").Value = 800
$ws.Range("D17").Value = 45013
$ws.Range("J17").Value = 1100
$ws.Range("K17").Value = 2000
$ws.Range("L17").Value = 2500
$ws.Range("M17").Value = 2250
$ws.Range("P17").Value = 750
$ws.Range("D18").Value = 44911
$ws.Range("J18").Value = 700
$ws.Range("K18").Value = 1800
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = 1900
$ws.Range("P18").Value = 633
$ws.Range("D19").Value = 44970
$ws.Range("J19").Value = 800
$ws.Range("D20").Value = 45006
$ws.Range("J20").Value = 1100
$ws.Range("D21").Value = 44827
$ws.Range("D22").Value = 44910
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 1800
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = 1900
$ws.Range("P22").Value = 633
$ws.Range("D23").Value = 45035
$ws.Range("J23").Value = 1100
$ws.Range("D24").Value = 45028
$ws.Range("J24").Value = 1000
$ws.Range("D25").Value = 44978
$ws.Range("J25").Value = 1000
$ws.Range("K25").Value = 1800
$ws.Range("L25").Value = 2000
$ws.Range("M25").Value = 1900
$ws.Range("P25").Value = 633
$ws.Range("D26").Value = 45020
$ws.Range("J26").Value = 1200
$ws.Range("D27").Value = 44992
$ws.Range("J27").Value = 1040
$ws.Range("D28").Value = 45034
$ws.Range("J28").Value = 1100
$ws.Range("K28").Value = 2000
$ws.Range("L28").Value = 2500
$ws.Range("M28").Value = 2250
$ws.Range("P28").Value = 750
$ws.Range("D29").Value = 44848
$ws.Range("J29").Value = 1000
$ws.Range("K29").Value = 1500
$ws.Range("L29").Value = 2000
$ws.Range("M29").Value = 1750
$ws.Range("P29").Value = 583
$ws.Range("D30").Value = 45007
$ws.Range("J30").Value = 1160
$ws.Range("D31").Value = 44685
$ws.Range("J31").Value = 400
$ws.Range("K31").Value = 1500
$ws.Range("L31").Value = 2000
$ws.Range("M31").Value = 1750
$ws.Range("P31").Value = 583
